$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 142.32259
$ws.Range("I33").Value = 85
$ws.Range("J33").Value = 440.4
$ws.Range("K33").Value = 85
$ws.Range("L33").Value = 440.4
$ws.Range("M33").Value = 144
$ws.Range("N33").Value = -898.4

$ws.Range("H137").Value = 2728.6309
$ws.Range("I137").Value = 759.5357
$ws.Range("J137").Value = 3713.1785
$ws.Range("K137").Value = 2278.6071
$ws.Range("L137").Value = 11139.5355
$ws.Range("M137").Value = 271.3928999999998
$ws.Range("N137").Value = -16239.5355

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2534.03
$ws.Range("I32").Value = 2418.2122
$ws.Range("J32").Value = 14000
$ws.Range("K32").Value = 2418.2122
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = -2131.2122
$ws.Range("N32").Value = -14574

$ws.Range("H45").Value = 2131.1538
$ws.Range("I45").Value = 1973.875
$ws.Range("J45").Value = 2382.8
$ws.Range("K45").Value = 1973.875
$ws.Range("L45").Value = 2382.8
$ws.Range("M45").Value = -1596.875
$ws.Range("N45").Value = -3136.8

$ws.Range("H61").Value = 3273.2068
$ws.Range("I61").Value = 1665.3077
$ws.Range("K61").Value = 1665.3077
$ws.Range("M61").Value = -1453.3077

$ws.Range("H74").Value = 2636.5881
$ws.Range("I74").Value = 2617.7693
$ws.Range("K74").Value = 2617.7693
$ws.Range("M74").Value = -1743.7693

$ws.Range("H77").Value = 2636.5881
$ws.Range("I77").Value = 2617.7693
$ws.Range("K77").Value = 13088.8465
$ws.Range("M77").Value = -8720.8465

$ws.Range("H110").Value = 1833.5
$ws.Range("I110").Value = 1275.8334
$ws.Range("K110").Value = 1275.8334
$ws.Range("M110").Value = 769.1666

$ws.Range("H122").Value = 1905.6364
$ws.Range("I122").Value = 1976.2
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 5928.6
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -3478.6
$ws.Range("N122").Value = -8500

$ws.Range("H132").Value = 7814241
$ws.Range("I132").Value = 11629036
$ws.Range("K132").Value = 34887108
$ws.Range("M132").Value = -34884578

$ws.Range("H136").Value = 3273.2068
$ws.Range("I136").Value = 1665.3077
$ws.Range("K136").Value = 4995.9231
$ws.Range("M136").Value = -2445.9231

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2002.0834
$ws.Range("I86").Value = 2257.1428
$ws.Range("J86").Value = 1645
$ws.Range("K86").Value = 2257.1428
$ws.Range("L86").Value = 1645
$ws.Range("M86").Value = -1134.1428
$ws.Range("N86").Value = -3891

$ws.Range("H89").Value = 2002.0834
$ws.Range("I89").Value = 2257.1428
$ws.Range("J89").Value = 1645
$ws.Range("K89").Value = 11285.714
$ws.Range("L89").Value = 8225
$ws.Range("M89").Value = -5669.714
$ws.Range("N89").Value = -19457

$ws.Range("H99").Value = 2044
$ws.Range("I99").Value = 2021.1111
$ws.Range("K99").Value = 2021.1111
$ws.Range("M99").Value = -523.1111000000001

$ws.Range("H107").Value = 1538.9429
$ws.Range("I107").Value = 1337.9
$ws.Range("J107").Value = 2745.2
$ws.Range("K107").Value = 1337.9
$ws.Range("L107").Value = 2745.2
$ws.Range("M107").Value = 582.0999999999999
$ws.Range("N107").Value = -6585.2

$ws.Range("H134").Value = 3938.5173
$ws.Range("I134").Value = 3119
$ws.Range("J134").Value = 4337.769
$ws.Range("K134").Value = 9357
$ws.Range("L134").Value = 13013.307
$ws.Range("M134").Value = -6822
$ws.Range("N134").Value = -18083.307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4305.8374
$ws.Range("I58").Value = 5105.552
$ws.Range("K58").Value = 5105.552
$ws.Range("M58").Value = -4902.552

$ws.Range("H64").Value = 30282.715
$ws.Range("J64").Value = 30282.715
$ws.Range("L64").Value = 30282.715
$ws.Range("N64").Value = -30778.715

$ws.Range("H67").Value = 30282.715
$ws.Range("J67").Value = 30282.715
$ws.Range("L67").Value = 30282.715
$ws.Range("N67").Value = -31998.715

$ws.Range("H99").Value = 2330.1667
$ws.Range("I99").Value = 2170.6667
$ws.Range("J99").Value = 2383.3333
$ws.Range("K99").Value = 2170.6667
$ws.Range("L99").Value = 2383.3333
$ws.Range("M99").Value = -672.6667000000002
$ws.Range("N99").Value = -5379.3333

$ws.Range("H105").Value = 2720.4285
$ws.Range("I105").Value = 3103
$ws.Range("J105").Value = 1764
$ws.Range("K105").Value = 3103
$ws.Range("L105").Value = 1764
$ws.Range("M105").Value = -1356
$ws.Range("N105").Value = -5258

$ws.Range("H122").Value = 174845.14
$ws.Range("J122").Value = 10457
$ws.Range("L122").Value = 31371
$ws.Range("N122").Value = -36271

$ws.Range("H126").Value = 2330.1667
$ws.Range("I126").Value = 2170.6667
$ws.Range("J126").Value = 2383.3333
$ws.Range("K126").Value = 6512.000100000001
$ws.Range("L126").Value = 7149.999899999999
$ws.Range("M126").Value = -4042.000100000001
$ws.Range("N126").Value = -12089.9999

$ws.Range("H132").Value = 42794.973
$ws.Range("I132").Value = 1663.1904
$ws.Range("K132").Value = 4989.5712
$ws.Range("M132").Value = -2459.5712

$ws.Range("H134").Value = 2156.1765
$ws.Range("I134").Value = 1006.4286
$ws.Range("K134").Value = 3019.2858
$ws.Range("M134").Value = -484.2857999999997

$ws.Range("H136").Value = 4305.8374
$ws.Range("I136").Value = 5105.552
$ws.Range("K136").Value = 15316.656
$ws.Range("M136").Value = -12766.656

$ws.Range("H138").Value = 45719.8
$ws.Range("J138").Value = 45719.8
$ws.Range("L138").Value = 45719.8
$ws.Range("N138").Value = -55999.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 7010.136
$ws.Range("J88").Value = 7010.136
$ws.Range("L88").Value = 21030.408
$ws.Range("N88").Value = -21886.408

$ws.Range("H91").Value = 7010.136
$ws.Range("J91").Value = 7010.136
$ws.Range("L91").Value = 21030.408
$ws.Range("N91").Value = -23994.408

$ws.Range("H113").Value = 4969.5654
$ws.Range("I113").Value = 6465.5884
$ws.Range("K113").Value = 19396.7652
$ws.Range("M113").Value = -17226.7652

$ws.Range("H121").Value = 163729.73
$ws.Range("I121").Value = 328
$ws.Range("J121").Value = 222087.5
$ws.Range("K121").Value = 984
$ws.Range("L121").Value = 666262.5
$ws.Range("M121").Value = 326
$ws.Range("N121").Value = -668882.5

$ws.Range("H137").Value = 33340650
$ws.Range("I137").Value = 2652.4375
$ws.Range("J137").Value = 71441220
$ws.Range("K137").Value = 7957.3125
$ws.Range("L137").Value = 214323660
$ws.Range("M137").Value = -2857.3125
$ws.Range("N137").Value = -214333860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 6067.696
$ws.Range("I113").Value = 7372.4116
$ws.Range("K113").Value = 7372.4116
$ws.Range("M113").Value = -5202.4116

$ws.Range("H114").Value = 42021.5
$ws.Range("J114").Value = 42021.5
$ws.Range("L114").Value = 42021.5
$ws.Range("N114").Value = -50699.5

$ws.Range("H135").Value = 34740
$ws.Range("J135").Value = 34740
$ws.Range("L135").Value = 34740
$ws.Range("N135").Value = -44880

$ws.Range("H138").Value = 54210
$ws.Range("J138").Value = 54210
$ws.Range("L138").Value = 54210
$ws.Range("N138").Value = -64490

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 24245
$ws.Range("J45").Value = 24490
$ws.Range("L45").Value = 24490
$ws.Range("N45").Value = -25304

$ws.Range("H112").Value = 43594.668
$ws.Range("J112").Value = 43594.668
$ws.Range("L112").Value = 43594.668
$ws.Range("N112").Value = -46548.668

$ws.Range("H136").Value = 2901.7307
$ws.Range("I136").Value = 2191.647
$ws.Range("J136").Value = 4243
$ws.Range("K136").Value = 6574.941
$ws.Range("L136").Value = 12729
$ws.Range("M136").Value = -4024.941
$ws.Range("N136").Value = -17829

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 279.6154
$ws.Range("I113").Value = 279.6154
$ws.Range("K113").Value = 838.8462000000001
$ws.Range("M113").Value = 1331.1538

$ws.Range("H126").Value = 1550663.2
$ws.Range("I126").Value = 1733036.6
$ws.Range("J126").Value = 490
$ws.Range("K126").Value = 5199109.800000001
$ws.Range("L126").Value = 1470
$ws.Range("M126").Value = -5196639.800000001
$ws.Range("N126").Value = -6410

$ws.Range("H132").Value = 1531.7826
$ws.Range("I132").Value = 1104.7693
$ws.Range("J132").Value = 2086.9
$ws.Range("K132").Value = 3314.3079
$ws.Range("L132").Value = 6260.700000000001
$ws.Range("M132").Value = -784.3078999999998
$ws.Range("N132").Value = -11320.7

$ws.Range("H136").Value = 18232.553
$ws.Range("I136").Value = 41671.52
$ws.Range("K136").Value = 125014.56
$ws.Range("M136").Value = -122464.56
